$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "marker" column (E), pushing
# the old E..I columns to G..K.
$ws.Columns("E:F").Insert()

# New column widths: new F ("functional_association") gets an explicit
# width; new E ("fresh_frozen") keeps the default width.
$ws.Columns("F").ColumnWidth = 24.95

$enhancers = [char]0xFEFF + "active promoters & enhancers"

# Populate the new "functional_association" column (F) grouped by assay,
# atac rows first, then H3K4me3, then H3K27ac - this matches the order in
# which the values were first authored.
$ws.Range("F9").Value = "open chromatin"
$ws.Range("F10").Value = "open chromatin"
$ws.Range("F11").Value = "open chromatin"
$ws.Range("F12").Value = "open chromatin"
$ws.Range("F13").Value = "open chromatin"
$ws.Range("F14").Value = "open chromatin"
$ws.Range("F15").Value = "open chromatin"

$ws.Range("F16").Value = "active promoters"
$ws.Range("F17").Value = "active promoters"
$ws.Range("F18").Value = "active promoters"
$ws.Range("F19").Value = "active promoters"

$ws.Range("F2").Value = $enhancers
$ws.Range("F3").Value = $enhancers
$ws.Range("F4").Value = $enhancers
$ws.Range("F5").Value = $enhancers
$ws.Range("F6").Value = $enhancers
$ws.Range("F7").Value = $enhancers
$ws.Range("F8").Value = $enhancers

$ws.Range("F1").Value = "functional_association"

# Populate the new "fresh_frozen" column (E).
$ws.Range("E1").Value = "fresh_frozen"

$ws.Range("E2").Value = "fresh"
$ws.Range("E3").Value = "fresh"
$ws.Range("E4").Value = "frozen"
$ws.Range("E5").Value = "frozen"
$ws.Range("E6").Value = "frozen"
$ws.Range("E7").Value = "frozen"
$ws.Range("E8").Value = "frozen"
$ws.Range("E9").Value = "fresh"
$ws.Range("E10").Value = "fresh"
$ws.Range("E11").Value = "frozen"
$ws.Range("E12").Value = "frozen"
$ws.Range("E13").Value = "frozen"
$ws.Range("E14").Value = "frozen"
$ws.Range("E15").Value = "frozen"
$ws.Range("E16").Value = "frozen"
$ws.Range("E17").Value = "frozen"
$ws.Range("E18").Value = "frozen"
$ws.Range("E19").Value = "frozen"

# Rows 16-19 use a slightly different (black, non-themed) font color on
# the new "fresh_frozen" cells.
$ws.Range("E16").Font.Color = 0
$ws.Range("E17").Font.Color = 0
$ws.Range("E18").Font.Color = 0
$ws.Range("E19").Font.Color = 0

# Update the view: scroll back to the top-left and select F25, matching
# the author's final selection/scroll position.
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("F25").Select()
